$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the first data row (row 2), pushing the
# existing data rows (3-10) down to (5-12).
$ws.Rows.Item(3).Resize(2).Insert()

# Copy the date-cell style (s="2") from D2 so the new D3/D4 cells keep the
# same date number format.
$ws.Range("D2").Copy()
$ws.Range("D3:D4").PasteSpecial(-4122)  # xlPasteFormats

# Row 3 - new "Primera" record for 2023-08-07
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 45145
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112039
$ws.Range("G3").Value = "Ciboulette"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2500
$ws.Range("N3").Value = "$/docena de atados"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 833
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = "Hortaliza"

# Row 4 - new "Segunda" record for 2023-08-07
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45145
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112039
$ws.Range("G4").Value = "Ciboulette"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 2000
$ws.Range("N4").Value = "$/docena de atados"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 667
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = "Hortaliza"
